$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet contains a product list (rows 4-9) sorted alphabetically by name,
# followed by a totals row and a footer row. This update adds three new
# products to the list (keeping alphabetical order) and refreshes the
# timestamp / totals accordingly.
#
# Current layout (rows 4-11):
#   4 BEBY RELIEF 25 MG  SUPP
#   5 DANSET 8MG/4ML 3 AMP.
#   6 جهاز محلول
#   7 سرنجات 5 سم
#   8 كالونا
#   9 محلول ملح
#   10 (totals row)
#   11 (footer row)
#
# New layout (rows 4-14):
#   4  BEBY RELIEF 25 MG  SUPP
#   5  CYMBATEX 20 MG 30 CAPS.   <- NEW
#   6  DANSET 8MG/4ML 3 AMP.
#   7  PANTOLOC 40MG 14 TAB      <- NEW
#   8  TRIACTIN 4MG 20 TAB       <- NEW
#   9  جهاز محلول
#   10 سرنجات 5 سم
#   11 كالونا
#   12 محلول ملح
#   13 (totals row)
#   14 (footer row)
# ---------------------------------------------------------------------------

# --- Insert "CYMBATEX 20 MG 30 CAPS." as new row 5 ------------------------
$ws.Range("A4:N4").Copy()
$ws.Rows.Item(5).Insert()
$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()

$ws.Range("B5").Value = "CYMBATEX 20 MG 30 CAPS."
$ws.Range("H5").Value = "0:2"
$ws.Range("L5").Value = 58
$ws.Range("N5").Value = "0:0"

# --- Insert "PANTOLOC 40MG 14 TAB" as new row 7 ----------------------------
# (row 6 is now DANSET, so the new product goes right after it)
$ws.Range("A6:N6").Copy()
$ws.Rows.Item(7).Insert()
$ws.Range("B7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()

$ws.Range("B7").Value = "PANTOLOC 40MG 14 TAB"
$ws.Range("H7").Value = "1:0"
$ws.Range("L7").Value = 51
$ws.Range("N7").Value = "0:2"

# --- Insert "TRIACTIN 4MG 20 TAB" as new row 8 -----------------------------
$ws.Range("A7:N7").Copy()
$ws.Rows.Item(8).Insert()
$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()

$ws.Range("B8").Value = "TRIACTIN 4MG 20 TAB"
$ws.Range("H8").Value = "1:1"
$ws.Range("L8").Value = 23
$ws.Range("N8").Value = "0:2"

# --- Renumber the "م" column (A4:A12) --------------------------------------
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9

# --- Row heights -------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 25.5
$ws.Rows.Item(14).RowHeight = 17.25

# --- Totals row (now row 13): recompute the sum of the price column -------
$ws.Range("K13").Value = 331

# --- Refresh the generated timestamp shown in the footer (row 14) ---------
$ws.Range("A14").Value = "Saturday, 17 January, 2026 9:31 AM"
